# Add new "power meter" global variables to the "IO Mapping" sheet.
#
# 1) Insert 4 rows at row 133 (Modbus MB_power_* mirror registers),
#    pushing the existing rows 133-208 down to 137-212.
# 2) Append 4 rows at the end (213-216) for the FEEDBACK_power_meter_*
#    readings (D4710-D4725).
# 3) Activate the "IO Mapping" sheet (it becomes the active tab), which
#    mirrors the tabSelected move from "Pump" -> "IO Mapping" and the
#    workbook-level activeTab bump from 4 -> 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IO Mapping")

# --- 1) Insert 4 new rows at 133, shifting everything below down ------
$ws.Range("A133:A136").EntireRow.Insert()

$newRows = @(
    @{ Row = 133; A = "MB_power_Ia"; B = "D25715" },
    @{ Row = 134; A = "MB_power_Ib"; B = "D25717" },
    @{ Row = 135; A = "MB_power_Ic"; B = "D25719" },
    @{ Row = 136; A = "MB_power_kw"; B = "D25729" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = "REAL"
    $ws.Cells.Item($row, 4).Value = -1
    $ws.Cells.Item($row, 5).Value = "x"

    # Column D on these rows keeps the quote-prefixed/centered style (s="4")
    # that the row above (132) uses - copy formats only, after the value is
    # already in place, so the numeric -1 is not converted back to text.
    $ws.Cells.Item(132, 4).Copy()
    $ws.Cells.Item($row, 4).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# --- 2) Append 4 rows (213-216) with the new FEEDBACK_power_meter_* regs
$feedbackRows = @(
    @{ Row = 213; A = "FEEDBACK_power_meter_KW"; B = "D4710"; F = "Total Power reading in KW " },
    @{ Row = 214; A = "FEEDBACK_power_meter_Ia"; B = "D4715"; F = "Total Power reading in Ia" },
    @{ Row = 215; A = "FEEDBACK_power_meter_Ib"; B = "D4720"; F = "Total Power reading in Ib" },
    @{ Row = 216; A = "FEEDBACK_power_meter_Ic"; B = "D4725"; F = "Total Power reading in Ic" }
)

foreach ($r in $feedbackRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = "REAL"
    $ws.Cells.Item($row, 3).HorizontalAlignment = -4131
    $ws.Cells.Item($row, 4).Value = 0
    $ws.Cells.Item($row, 5).Value = "x"
    $ws.Cells.Item($row, 6).Value = $r.F
}

# --- 3) Make "IO Mapping" the active/selected tab, matching the view state
$ws.Cells.Item(128, 1).Select()
$ws.Application.ActiveWindow.ScrollRow = 128
$ws.Range("B137").Select()
$ws.Activate()
